$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "79.421.16"
$ws.Range("E2").Value = "  +4.02%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.190.40"
$ws.Range("E3").Value = "  +5.00%  "
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "204.15"
$ws.Range("E5").Value = "  +1.81%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "636.29"
$ws.Range("E6").Value = "  +2.11%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.245"
$ws.Range("E8").Value = "  +19.40%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.604"
$ws.Range("E9").Value = "  +9.73%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.189.21"
$ws.Range("E10").Value = "  +5.02%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.645"
$ws.Range("E11").Value = "  +46.51%  "
$ws.Range("E12").Value = "  +2.96%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.0000239"
$ws.Range("E13").Value = "  +24.46%  "
$ws.Range("B14").Value = "Toncoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.44"
$ws.Range("E14").Value = "  +3.85%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "3.768.90"
$ws.Range("E15").Value = "  +4.71%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "32.11"
$ws.Range("E16").Value = "  +10.00%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "79.250.36"
$ws.Range("E17").Value = "  +3.82%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.174.54"
$ws.Range("E18").Value = "  +4.07%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "14.63"
$ws.Range("E19").Value = "  +8.05%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "3.00"
$ws.Range("E20").Value = "  +30.58%  "
$ws.Range("E21").Value = "  +3.46%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "432.44"
$ws.Range("E22").Value = "  +15.32%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.24"
$ws.Range("E23").Value = "  +20.70%  "
$ws.Range("B24").Value = "NEARProtocol"
$ws.Range("C24").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "4.82"
$ws.Range("E24").Value = "  +10.24%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "3.358.81"
$ws.Range("E25").Value = "  +5.00%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "76.99"
$ws.Range("E26").Value = "  +4.76%  "
$ws.Range("B27").Value = "Aptos"
$ws.Range("C27").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.09"
$ws.Range("E27").Value = "  +13.01%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.999"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("E29").Value = "  +9.63%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "9.09"
$ws.Range("E30").Value = "  +10.36%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  +5.02%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "523.46"
$ws.Range("E33").Value = "  +5.45%  "
$ws.Range("E34").Value = "  +2.46%  "
$ws.Range("E35").Value = "  +22.21%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "23.30"
$ws.Range("E36").Value = "  +12.91%  "
$ws.Range("E37").Value = "  +16.25%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.999"
$ws.Range("E38").Value = "  -0.04%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.406"
$ws.Range("E39").Value = "  +5.82%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "165.15"
$ws.Range("E40").Value = "  +1.42%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "20.02"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("B42").Value = "USDe"
$ws.Range("C42").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "191.41"
$ws.Range("E43").Value = "  +0.95%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.51"
$ws.Range("E44").Value = "  +7.87%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.803"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("E46").Value = "  +8.14%  "
$ws.Range("E47").Value = "  +4.64%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "43.40"
$ws.Range("E48").Value = "  +3.32%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "25.70"
$ws.Range("E49").Value = "  +15.12%  "
$ws.Range("B50").Value = "ARBITRUM"
$ws.Range("C50").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.637"
$ws.Range("E50").Value = "  +4.74%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.53"
$ws.Range("E51").Value = "  +2.51%  "
